# Updated symbol list on Mon Jan 30 15:27:40 UTC 2023 with GitHub Actions
# Applies updated Price (D) and Volume(1h) (E) values to Sheet1 rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Leading apostrophe forces Excel to store the numeric/percent-looking
    # string as literal text instead of coercing it to a number.
    $range.Value = "'" + $text
    # Re-apply the Normal style so the quote-prefix flag picked up above
    # does not leave a stray number-format/style behind on the cell.
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "309.35"
Set-TextValue $ws.Range("E2") "-3.19%"
Set-TextValue $ws.Range("D3") "37.89"
Set-TextValue $ws.Range("E3") "-4.29%"
Set-TextValue $ws.Range("D4") "5.084"
Set-TextValue $ws.Range("E4") "-0.65%"
Set-TextValue $ws.Range("D5") "0.07888"
Set-TextValue $ws.Range("E5") "-3.65%"
Set-TextValue $ws.Range("D6") "1.977"
Set-TextValue $ws.Range("E6") "-1.95%"
Set-TextValue $ws.Range("D7") "4.363"
Set-TextValue $ws.Range("E7") "2.09%"
Set-TextValue $ws.Range("D8") "8.293"
Set-TextValue $ws.Range("E8") "0.13%"
Set-TextValue $ws.Range("D9") "3.036"
Set-TextValue $ws.Range("E9") "-4.14%"
Set-TextValue $ws.Range("D10") "0.9314"
Set-TextValue $ws.Range("E10") "-0.32%"
Set-TextValue $ws.Range("E11") "-7.62%"
Set-TextValue $ws.Range("D12") "0.1971"
Set-TextValue $ws.Range("E12") "-1.14%"
Set-TextValue $ws.Range("D13") "0.08815"
Set-TextValue $ws.Range("E13") "-3.18%"
Set-TextValue $ws.Range("D14") "0.03424"
Set-TextValue $ws.Range("E14") "-4.31%"
Set-TextValue $ws.Range("D15") "0.09738"
Set-TextValue $ws.Range("E15") "-0.68%"
Set-TextValue $ws.Range("D16") "0.001388"
Set-TextValue $ws.Range("E16") "-1.21%"
Set-TextValue $ws.Range("D17") "0.005898"
Set-TextValue $ws.Range("E17") "-5.13%"
Set-TextValue $ws.Range("E18") "1,776.45%"
Set-TextValue $ws.Range("E19") "-1.93%"
Set-TextValue $ws.Range("D20") "0.3475"
Set-TextValue $ws.Range("E20") "0.46%"
Set-TextValue $ws.Range("D21") "0.1296"
Set-TextValue $ws.Range("E21") "0.78%"
Set-TextValue $ws.Range("D22") "5.005"
Set-TextValue $ws.Range("E22") "2.03%"
Set-TextValue $ws.Range("D23") "0.2487"
Set-TextValue $ws.Range("E23") "1.46%"
Set-TextValue $ws.Range("D24") "0.04314"
Set-TextValue $ws.Range("E24") "-0.47%"
Set-TextValue $ws.Range("D25") "0.001217"
Set-TextValue $ws.Range("E25") "-0.70%"
Set-TextValue $ws.Range("D26") "0.004616"
Set-TextValue $ws.Range("E26") "-3.54%"
Set-TextValue $ws.Range("D27") "0.0001351"
Set-TextValue $ws.Range("E27") "3.84%"
Set-TextValue $ws.Range("D39") "0.02280"
Set-TextValue $ws.Range("E39") "2.22%"
Set-TextValue $ws.Range("D40") "0.05026"
Set-TextValue $ws.Range("E40") "-4.58%"
Set-TextValue $ws.Range("D41") "0.007501"
Set-TextValue $ws.Range("E41") "0.01%"
Set-TextValue $ws.Range("D42") "0.009854"
Set-TextValue $ws.Range("E42") "-1.63%"
Set-TextValue $ws.Range("D43") "0.1359"
Set-TextValue $ws.Range("E43") "-1.43%"
Set-TextValue $ws.Range("D44") "0.002042"
Set-TextValue $ws.Range("E44") "-5.57%"
Set-TextValue $ws.Range("D45") "0.008789"
Set-TextValue $ws.Range("E45") "-11.11%"
Set-TextValue $ws.Range("D46") "0.00006590"
Set-TextValue $ws.Range("E46") "1.26%"
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "-0.14%"
Set-TextValue $ws.Range("D48") "0.003000"
Set-TextValue $ws.Range("D50") "0.00002099"
Set-TextValue $ws.Range("E50") "-0.14%"
Set-TextValue $ws.Range("D51") "0.0001999"
Set-TextValue $ws.Range("E51") "-0.14%"
